$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "파비클래스, PDSI를 거치지 않은 질문을 받지 말아야 할 이유"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/pabiiclass-pdsi-not-answering/#utm_source=rss&utm_medium=rss&utm_campaign=pabiiclass-pdsi-not-answering"

$ws.Range("D28").Value = "Multiobjective Optimization for stiffness and position control in a soft robot arm module"
$ws.Range("E28").Value = "https://ropiens.tistory.com/140"

$ws.Range("D46").Value = "[SK바이오팜] 2021년 07월, 생물정보학(Bioinformatics 채용), Digital Healthcare(Computer-Aided Drug Design, CADD) 경력 구성원 영입"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/412"
